# Update "想去人数" (F column) counts across the sheets to reflect the
# latest scrape snapshot (gh-pages output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- 展览 (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 27010
$ws1.Range("F5").Value  = 644
$ws1.Range("F7").Value  = 563
$ws1.Range("F8").Value  = 233
$ws1.Range("F10").Value = 474
$ws1.Range("F13").Value = 315
$ws1.Range("F15").Value = 476
$ws1.Range("F17").Value = 1616
$ws1.Range("F19").Value = 505
$ws1.Range("F20").Value = 137
$ws1.Range("F21").Value = 455
$ws1.Range("F22").Value = 8

# --- 演出 (sheet2) ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value  = 4522
$ws2.Range("F11").Value = 452

# --- 本地生活 (sheet3) ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value  = 269

# --- 全部类型 (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 269
$ws4.Range("F5").Value  = 27010
$ws4.Range("F6").Value  = 4522
$ws4.Range("F9").Value  = 644
$ws4.Range("F18").Value = 452
$ws4.Range("F19").Value = 563
$ws4.Range("F22").Value = 233
$ws4.Range("F24").Value = 474
$ws4.Range("F28").Value = 315
$ws4.Range("F32").Value = 476
$ws4.Range("F35").Value = 1616
$ws4.Range("F37").Value = 505
$ws4.Range("F39").Value = 137
$ws4.Range("F40").Value = 455
$ws4.Range("F41").Value = 8
